$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values parse as plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts the assigned string into
# a floating point number (losing the exact display text, e.g. trailing zeros).

$ws.Range("D2").Value = '30.483.56'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.889.38'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.73'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4711'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2898'
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06495'
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.21'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07752'
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '1.889.84'
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.73'
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7256'
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.189'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '281.14'
$ws.Range("E16").Value = '  +2.88%  '
$ws.Range("D17").Value = '30.470.77'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.06'
$ws.Range("E18").Value = '  -1.81%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007475'
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("D21").Value = '2.136.97'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.257'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.82'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.075'
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.88'
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.892'
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09719'
$ws.Range("E29").Value = '  -3.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.332'
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.467'
$ws.Range("E31").Value = '  -2.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.280'
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.148'
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04861'
$ws.Range("E34").Value = '  +1.30%  '
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6936'
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01885'
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.818'
$ws.Range("E39").Value = '  +2.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.49'
$ws.Range("E40").Value = '  +3.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.213'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.982'
$ws.Range("E42").Value = '  +0.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4248'
$ws.Range("E43").Value = '  +1.92%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8240'
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.34'
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.566'
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.961'
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.11'
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05748'
$ws.Range("E51").Value = '  +1.65%  '
